# Auto-generated edit script applying the diff changes to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F45").Value = 67
$ws.Range("G45").Value = 12923.63
$ws.Range("B71").Value = 56512.15
$ws.Range("F105").Value = 223
$ws.Range("G105").Value = 18100.91
$ws.Range("F136").Value = 11
$ws.Range("G136").Value = 2961.09
$ws.Range("B143").Value = 280431.19
$ws.Range("B213").Value = 53925
$ws.Range("B214").Value = 57756
$ws.Range("F235").Value = 35
$ws.Range("G235").Value = 2954.7
$ws.Range("F237").Value = 11
$ws.Range("G237").Value = 1224.63
$ws.Range("F239").Value = 67
$ws.Range("G239").Value = 5494
$ws.Range("B250").Value = 29657.64
$ws.Range("F257").Value = 25
$ws.Range("G257").Value = 1875.5
$ws.Range("B262").Value = 15787.34
$ws.Range("B387").Value = 61610
$ws.Range("D387").Value = 102.71
$ws.Range("E387").Value = 122.71
$ws.Range("F387").Value = 248
$ws.Range("G387").Value = 25472.08
$ws.Range("B388").Value = 57077
$ws.Range("D388").Value = 93.08
$ws.Range("E388").Value = 111.2
$ws.Range("F388").Value = 1
$ws.Range("G388").Value = 93.08
$ws.Range("F397").Value = 26
$ws.Range("G397").Value = 3136.9
$ws.Range("F405").Value = 113
$ws.Range("G405").Value = 15488.91
$ws.Range("F406").Value = 132
$ws.Range("G406").Value = 19044.96
$ws.Range("F411").Value = 151
$ws.Range("G411").Value = 7749.32
$ws.Range("F425").Value = 3
$ws.Range("G425").Value = 257.16
$ws.Range("F430").Value = 16
$ws.Range("G430").Value = 1665.44
$ws.Range("F434").Value = 175
$ws.Range("G434").Value = 17692.5
$ws.Range("F442").Value = 170
$ws.Range("G442").Value = 10052.1
$ws.Range("F447").Value = 14
$ws.Range("G447").Value = 303.66
$ws.Range("B467").Value = 422313.91
$ws.Range("F536").Value = 105
$ws.Range("G536").Value = 10143
$ws.Range("B542").Value = 56840.67
$ws.Range("F563").Value = 76
$ws.Range("G563").Value = 14169.44
$ws.Range("B571").Value = 55305.68
$ws.Range("F574").Value = 12
$ws.Range("G574").Value = 606.6
$ws.Range("B582").Value = 1702.3
$ws.Range("F672").Value = 666
$ws.Range("G672").Value = 13220.1
$ws.Range("F673").Value = 473
$ws.Range("G673").Value = 3169.1
$ws.Range("B677").Value = 43512.75
$ws.Range("F689").Value = 88
$ws.Range("G689").Value = 5447.2
$ws.Range("B701").Value = 32966.18
$ws.Range("F725").Value = 131
$ws.Range("G725").Value = 5939.54
$ws.Range("F737").Value = 27
$ws.Range("G737").Value = 117.45
$ws.Range("B739").Value = 8707.809999999999
$ws.Range("F763").Value = 127
$ws.Range("G763").Value = 16579.85
$ws.Range("F767").Value = 112
$ws.Range("G767").Value = 3046.4
$ws.Range("F768").Value = 112
$ws.Range("G768").Value = 3046.4
$ws.Range("F769").Value = 46
$ws.Range("G769").Value = 1251.2
$ws.Range("B770").Value = 70304.98
$ws.Range("F791").Value = 3
$ws.Range("G791").Value = 3651.99
$ws.Range("B792").Value = 84896.02
$ws.Range("F802").Value = 185
$ws.Range("G802").Value = 7988.3
$ws.Range("B804").Value = 37866.91
$ws.Range("F833").Value = 12
$ws.Range("G833").Value = 1101
$ws.Range("B839").Value = 5017.1
$ws.Range("F863").Value = 21
$ws.Range("G863").Value = 1841.7
$ws.Range("B866").Value = 4938.42
$ws.Range("F869").Value = 21
$ws.Range("G869").Value = 4950.54
$ws.Range("F877").Value = 19
$ws.Range("G877").Value = 4663.93
$ws.Range("B880").Value = 29579.22
$ws.Range("F882").Value = 0
$ws.Range("G882").Value = 0
$ws.Range("F884").Value = 0
$ws.Range("G884").Value = 0
$ws.Range("F889").Value = 0
$ws.Range("G889").Value = 0
$ws.Range("F891").Value = 222
$ws.Range("G891").Value = 15442.32
$ws.Range("F893").Value = 94
$ws.Range("G893").Value = 3896.3
$ws.Range("F895").Value = 9
$ws.Range("G895").Value = 1360.17
$ws.Range("F898").Value = 0
$ws.Range("G898").Value = 0
$ws.Range("F899").Value = 2
$ws.Range("G899").Value = 241.42
$ws.Range("B901").Value = 52865.63
$ws.Range("F910").Value = 4
$ws.Range("G910").Value = 340.8
$ws.Range("F920").Value = 246
$ws.Range("G920").Value = 7409.52
$ws.Range("F926").Value = 40
$ws.Range("G926").Value = 1473.2
$ws.Range("B933").Value = 39779.61
$ws.Range("F935").Value = 3
$ws.Range("G935").Value = 322.74
$ws.Range("F936").Value = 73
$ws.Range("G936").Value = 2730.2
$ws.Range("F939").Value = 156
$ws.Range("G939").Value = 5834.4
$ws.Range("B942").Value = 14784.52
$ws.Range("F944").Value = 1
$ws.Range("G944").Value = 522.85
$ws.Range("B948").Value = 12691.45
$ws.Range("F961").Value = 3
$ws.Range("G961").Value = 5241.48
$ws.Range("F983").Value = 3
$ws.Range("G983").Value = 2911.77
$ws.Range("B984").Value = 128438.19
$ws.Range("F994").Value = 0
$ws.Range("G994").Value = 0
$ws.Range("B997").Value = 3591.88
$ws.Range("F999").Value = 1459
$ws.Range("G999").Value = 237977.49
$ws.Range("F1003").Value = 190
$ws.Range("G1003").Value = 12825
$ws.Range("B1005").Value = 280437.77
$ws.Range("B1012").Value = 2692249.17
$ws.Range("B1013").Value = 2692249.17

Write-Output "Applied 146 cell updates"
